$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 with data, mirroring formatting of existing data rows (2 and 3)
$ws.Cells.Item(4, 1).Value = 10
$ws.Cells.Item(4, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value = "La Araucanía"

# Column D uses the same style/number format as D2/D3 (date)
$ws.Cells.Item(4, 4).Value = 44452
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

$ws.Cells.Item(4, 5).Value = 9
$ws.Cells.Item(4, 6).Value = 100112042
$ws.Cells.Item(4, 7).Value = "Locoto"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 120
$ws.Cells.Item(4, 11).Value = 2300
$ws.Cells.Item(4, 12).Value = 2300
$ws.Cells.Item(4, 13).Value = 2300
$ws.Cells.Item(4, 14).Value = "$/kilo"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 2300
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"
